# Applies a row-content permutation to rows 2-13 of the active sheet.
# Only the columns that actually differ between the observation rows
# (A, B, E, F, G, H, Q, R, AC) are re-shuffled; row 3 is left untouched
# since its content does not move.
#
# Mapping: new row -> old row that its content comes from.
#   2<-9  3<-3  4<-12  5<-10  6<-13  7<-11  8<-7  9<-2  10<-8  11<-5  12<-4  13<-6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","E","F","G","H","Q","R","AC")

# Snapshot current values of the varying columns for every data row.
$snapshot = @{}
foreach ($r in 2..13) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowData
}

# Destination row -> source row (content moves from source into destination).
$mapping = @{
    2  = 9
    3  = 3
    4  = 12
    5  = 10
    6  = 13
    7  = 11
    8  = 7
    9  = 2
    10 = 8
    11 = 5
    12 = 4
    13 = 6
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcData = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $value = $srcData[$col]
        if ($col -eq "AC") {
            if ($null -eq $value -or $value -eq "") {
                $ws.Range("AC$destRow").ClearContents()
            } else {
                $ws.Range("AC$destRow").Value = $value
            }
        } else {
            $ws.Range("$col$destRow").Value = $value
        }
    }
}
